# Generate Report for Handoff
# Updates localization status cells from "In Translation" to "Ready for handoff"
# and refreshes the handoff timestamps, then widens the affected "Status"
# columns to fit the new (longer) text.

$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Worksheets.Item("Overview")
$wsZhCn     = $wb.Worksheets.Item("zh-cn")
$wsDeDe     = $wb.Worksheets.Item("de-de")

# --- Status text: "In Translation" -> "Ready for handoff" ---
$wsOverview.Range("E2").Value = "Ready for handoff"
$wsOverview.Range("F2").Value = "Ready for handoff"
$wsZhCn.Range("C2").Value = "Ready for handoff"
$wsDeDe.Range("C2").Value = "Ready for handoff"

# --- Handoff timestamps refreshed to reflect the new handoff generation ---
$wsOverview.Range("G2").Value = "2016-08-24 14:43:47"
$wsDeDe.Range("H2").Value = "2016-08-24 14:43:47"
$wsZhCn.Range("H2").Value = "2016-08-24 14:43:43"

# --- Widen the "Status" columns so the longer text fits ---
$wsOverview.Columns.Item(5).ColumnWidth = 16.33
$wsOverview.Columns.Item(6).ColumnWidth = 16.33
$wsZhCn.Columns.Item(3).ColumnWidth = 16.33
$wsDeDe.Columns.Item(3).ColumnWidth = 16.33
